$d = $word.ActiveDocument

# The document starts as a single paragraph:
#   "File is added and time is @ 4:11" [_GoBack bookmark]
# It needs to become two paragraphs:
#   "File is added and time is @ 4:11"
#   "Updated "@4:21[_GoBack]" "   (curly quotes; bookmark collapsed mid-run)

# Split the existing sentence into its own paragraph by replacing it with
# itself plus a paragraph break. Word keeps the (hidden) "_GoBack"
# bookmark collapsed at the point of the last edit, which lands right at
# the new paragraph boundary.
$find = $d.Content.Find
$find.Execute("File is added and time is @ 4:11", $true, $false, $false, $false, $false, $true, 1, $false, "File is added and time is @ 4:11^p", 2) | Out-Null

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

# Append the closing-quote/space chunk at the end of the (now empty)
# second paragraph first ...
$tail = $p2.Range
$tail.InsertAfter([char]0x201D + " ")

# ... then insert the "Updated "@4:21" lead-in right before it, at the
# paragraph boundary (the former end of paragraph 1). Using
# InsertBefore/InsertAfter (rather than assigning Range.Text) keeps each
# chunk as its own run instead of merging into one run that would later
# need to be re-split around the bookmark, so the "_GoBack" bookmark
# naturally ends up sitting, collapsed, between the two runs.
$head = $d.Range($p1.Range.End, $p1.Range.End)
$head.InsertBefore("Updated " + [char]0x201C + "@4:21")
